$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format from C9 onto the new date cells C10, C11, C13
$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C13").PasteSpecial(-4122)  # xlPasteFormats

# Set the new date values
$ws.Range("C11").Value2 = 45241
$ws.Range("C13").Value2 = 45245

# Set the new text values - write E13 before E11 so the shared-string
# indices come out in the expected order (2 = Finition..., 3 = Mockup Fini)
$ws.Range("E13").Value = "Finition de la partie intro du CDC"
$ws.Range("E11").Value = "Mockup Fini"

# Update the active selection to match the edited workbook
$ws.Range("E18").Select()
